# Applies the crypto price/volume refresh captured in the commit
# "Updated cryptos list ... with GitHub Actions".
# For each changed row we update Price (D) and Volume(1h) (E); rows 27/28
# additionally swap Coin (B) and Link (C) because Stellar/Cosmos traded places.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '29.258.82'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '1.846.05'
$ws.Range("E3").Value = '  +0.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9990'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.01'
$ws.Range("E5").Value = '  +0.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6628'
$ws.Range("E6").Value = '  -1.00%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.98'
$ws.Range("E8").Value = '  +6.90%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07448'
$ws.Range("E9").Value = '  +0.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2959'
$ws.Range("E10").Value = '  -0.08%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '23.36'
$ws.Range("E11").Value = '  +2.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07766'
$ws.Range("E12").Value = '  +0.62%  '

$ws.Range("D13").Value = '1.848.21'
$ws.Range("E13").Value = '  -4.00%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.026'
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6737'
$ws.Range("E15").Value = '  -0.69%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.45'
$ws.Range("E16").Value = '  -3.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.192'
$ws.Range("E17").Value = '  -0.05%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008759'
$ws.Range("E18").Value = '  +5.98%  '

$ws.Range("D19").Value = '29.254.24'
$ws.Range("E19").Value = '  +2.08%  '

$ws.Range("D20").Value = '2.094.06'
$ws.Range("E20").Value = '  +12.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.56'
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.06'
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("E23").Value = '  +0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.173'
$ws.Range("E24").Value = '  -0.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.0000'
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '158.85'
$ws.Range("E26").Value = '  -0.93%  '

$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1407'
$ws.Range("E27").Value = '  -0.16%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.635'
$ws.Range("E28").Value = '  -0.88%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.514'
$ws.Range("E30").Value = '  +0.39%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.147'
$ws.Range("E31").Value = '  -1.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.067'

$ws.Range("E33").Value = '  +0.47%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05331'
$ws.Range("E34").Value = '  -0.70%  '

$ws.Range("E35").Value = '  -0.47%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7493'
$ws.Range("E36").Value = '  -0.61%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.159'
$ws.Range("E37").Value = '  +1.74%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.654'
$ws.Range("E38").Value = '  -0.80%  '

$ws.Range("D39").Value = '1.314.36'
$ws.Range("E39").Value = '  -1.12%  '

$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.757'
$ws.Range("E41").Value = '  +0.90%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.390'
$ws.Range("E42").Value = '  +6.86%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9068'
$ws.Range("E43").Value = '  -1.58%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9995'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '103.73'
$ws.Range("E45").Value = '  +0.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.07964'
$ws.Range("E46").Value = '  +1.95%  '

$ws.Range("D47").Value = '1.992.39'
$ws.Range("E47").Value = '  +8.92%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '65.35'
$ws.Range("E48").Value = '  +1.96%  '

$ws.Range("E49").Value = '  -0.73%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.5142'
$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.757'
$ws.Range("E51").Value = '  -0.42%  '
